$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 114066.555
$ws.Range("I64").Value = 1000000
$ws.Range("J64").Value = 3324.875
$ws.Range("K64").Value = 1000000
$ws.Range("L64").Value = 3324.875
$ws.Range("M64").Value = -999752
$ws.Range("N64").Value = -3820.875

$ws.Range("H67").Value = 114066.555
$ws.Range("I67").Value = 1000000
$ws.Range("J67").Value = 3324.875
$ws.Range("K67").Value = 1000000
$ws.Range("L67").Value = 3324.875
$ws.Range("M67").Value = -999142
$ws.Range("N67").Value = -5040.875

$ws.Range("H70").Value = 1100
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1100
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -3840

$ws.Range("H73").Value = 1100
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1100
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -5172

$ws.Range("H129").Value = 933698.7
$ws.Range("J129").Value = 1069575.4
$ws.Range("L129").Value = 3208726.2
$ws.Range("N129").Value = -3218726.2

$ws.Range("H137").Value = 1710.2222
$ws.Range("I137").Value = 1668
$ws.Range("J137").Value = 1820
$ws.Range("K137").Value = 5004
$ws.Range("L137").Value = 5460
$ws.Range("M137").Value = -2454
$ws.Range("N137").Value = -10560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1248.1351
$ws.Range("I61").Value = 1152.7
$ws.Range("J61").Value = 1657.1428
$ws.Range("K61").Value = 1152.7
$ws.Range("L61").Value = 1657.1428
$ws.Range("M61").Value = -940.7
$ws.Range("N61").Value = -2081.1428

$ws.Range("H74").Value = 2061.8857
$ws.Range("I74").Value = 959.8929000000001
$ws.Range("J74").Value = 6469.857
$ws.Range("K74").Value = 959.8929000000001
$ws.Range("L74").Value = 6469.857
$ws.Range("M74").Value = -85.89290000000005
$ws.Range("N74").Value = -8217.857

$ws.Range("H77").Value = 2061.8857
$ws.Range("I77").Value = 959.8929000000001
$ws.Range("J77").Value = 6469.857
$ws.Range("K77").Value = 4799.4645
$ws.Range("L77").Value = 32349.285
$ws.Range("M77").Value = -431.4645
$ws.Range("N77").Value = -41085.285

$ws.Range("H132").Value = 3126.8147
$ws.Range("I132").Value = 2734.0557
$ws.Range("K132").Value = 8202.167099999999
$ws.Range("M132").Value = -5672.167099999999

$ws.Range("H136").Value = 1248.1351
$ws.Range("I136").Value = 1152.7
$ws.Range("J136").Value = 1657.1428
$ws.Range("K136").Value = 3458.1
$ws.Range("L136").Value = 4971.428400000001
$ws.Range("M136").Value = -908.1000000000004
$ws.Range("N136").Value = -10071.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 502242.75
$ws.Range("I105").Value = 502490
$ws.Range("J105").Value = 501995.5
$ws.Range("K105").Value = 502490
$ws.Range("L105").Value = 501995.5
$ws.Range("M105").Value = -500743
$ws.Range("N105").Value = -505489.5

$ws.Range("H134").Value = 2906.05
$ws.Range("I134").Value = 2955.9443
$ws.Range("K134").Value = 8867.832900000001
$ws.Range("M134").Value = -6332.832900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27077.44
$ws.Range("I31").Value = 60430.766
$ws.Range("J31").Value = 3452.1667
$ws.Range("K31").Value = 60430.766
$ws.Range("L31").Value = 3452.1667
$ws.Range("M31").Value = -60135.766
$ws.Range("N31").Value = -4042.1667

$ws.Range("H34").Value = 27077.44
$ws.Range("I34").Value = 60430.766
$ws.Range("J34").Value = 3452.1667
$ws.Range("K34").Value = 60430.766
$ws.Range("L34").Value = 3452.1667
$ws.Range("M34").Value = -60228.766
$ws.Range("N34").Value = -3856.1667

$ws.Range("H58").Value = 8938
$ws.Range("I58").Value = 1177.2413
$ws.Range("J58").Value = 65203.5
$ws.Range("K58").Value = 1177.2413
$ws.Range("L58").Value = 65203.5
$ws.Range("M58").Value = -974.2412999999999
$ws.Range("N58").Value = -65609.5

$ws.Range("H132").Value = 4174.533
$ws.Range("I132").Value = 4264.048
$ws.Range("K132").Value = 12792.144
$ws.Range("M132").Value = -10262.144

$ws.Range("H134").Value = 1462.1724
$ws.Range("I134").Value = 1481.9546
$ws.Range("K134").Value = 4445.8638
$ws.Range("M134").Value = -1910.8638

$ws.Range("H136").Value = 8938
$ws.Range("I136").Value = 1177.2413
$ws.Range("J136").Value = 65203.5
$ws.Range("K136").Value = 3531.7239
$ws.Range("L136").Value = 195610.5
$ws.Range("M136").Value = -981.7239
$ws.Range("N136").Value = -200710.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 500
$ws.Range("I41").Value = 500
$ws.Range("K41").Value = 1500
$ws.Range("M41").Value = -1162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 64574.293
$ws.Range("I70").Value = 110251.9
$ws.Range("J70").Value = 6716
$ws.Range("K70").Value = 110251.9
$ws.Range("L70").Value = 6716
$ws.Range("M70").Value = -109981.9
$ws.Range("N70").Value = -7256

$ws.Range("H73").Value = 64574.293
$ws.Range("I73").Value = 110251.9
$ws.Range("J73").Value = 6716
$ws.Range("K73").Value = 110251.9
$ws.Range("L73").Value = 6716
$ws.Range("M73").Value = -109315.9
$ws.Range("N73").Value = -8588

$ws.Range("H132").Value = 3163.0476
$ws.Range("I132").Value = 2809.1333
$ws.Range("K132").Value = 8427.3999
$ws.Range("M132").Value = -5897.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3697.318
$ws.Range("I132").Value = 3754.5789
$ws.Range("K132").Value = 11263.7367
$ws.Range("M132").Value = -8733.736699999999

$ws.Range("H136").Value = 1554.6522
$ws.Range("I136").Value = 1408.4445
$ws.Range("K136").Value = 4225.333500000001
$ws.Range("M136").Value = -1675.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9919.200000000001
$ws.Range("J74").Value = 9919.200000000001
$ws.Range("L74").Value = 9919.200000000001
$ws.Range("N74").Value = -11791.2

$ws.Range("H77").Value = 9919.200000000001
$ws.Range("J77").Value = 9919.200000000001
$ws.Range("L77").Value = 29757.6
$ws.Range("N77").Value = -39117.60000000001

$ws.Range("H81").Value = 334529.84
$ws.Range("I81").Value = 500350
$ws.Range("J81").Value = 251619.75
$ws.Range("K81").Value = 1000700
$ws.Range("L81").Value = 503239.5
$ws.Range("M81").Value = -999639
$ws.Range("N81").Value = -505361.5

$ws.Range("H84").Value = 334529.84
$ws.Range("I84").Value = 500350
$ws.Range("J84").Value = 251619.75
$ws.Range("K84").Value = 5003500
$ws.Range("L84").Value = 2516197.5
$ws.Range("M84").Value = -4998196
$ws.Range("N84").Value = -2526805.5

$ws.Range("H122").Value = 2241.05
$ws.Range("I122").Value = 1264.2222
$ws.Range("K122").Value = 3792.6666
$ws.Range("M122").Value = -1342.6666

$ws.Range("H132").Value = 1376.8462
$ws.Range("I132").Value = 990.4545000000001
$ws.Range("K132").Value = 2971.3635
$ws.Range("M132").Value = -441.3635000000004

$ws.Range("H136").Value = 673.9545000000001
$ws.Range("I136").Value = 595.94446
$ws.Range("J136").Value = 1025
$ws.Range("K136").Value = 1787.83338
$ws.Range("L136").Value = 3075
$ws.Range("M136").Value = 762.16662
$ws.Range("N136").Value = -8175
